$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.753.79'
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").Value = '1.887.32'
$ws.Range("E3").Value = '  -1.01%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7524'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.23%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '239.24'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.51%  '
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3027'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.13'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06779'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.51%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07934'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.24%  '
$ws.Range("B12").Value = 'Polygon'
$ws.Range("C12").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7396'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.57%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.879.71'
$ws.Range("E13").Value = '  -0.60%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.130'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.22%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.13'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.25%  '
$ws.Range("D16").Value = '29.762.48'
$ws.Range("E16").Value = '  -1.17%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.82'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.95%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.897'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '241.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.18%  '
$ws.Range("E20").Value = '  -1.57%  '
$ws.Range("E21").Value = '  -0.25%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.868'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '165.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.49%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.165'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.97%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '18.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.47%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1269'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.003'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.66%  '
$ws.Range("E29").Value = '  +2.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.511'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.227'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.985'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05198'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.33%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.243'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.44%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7229'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.87%  '
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("E37").Value = '  -1.70%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.760'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.83%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.115'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.59%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.4366'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.08%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.99'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.80%  '
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.873'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.13%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8243'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '7.543'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '99.29'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.63%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.676'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.14%  '
$ws.Range("D48").Value = '2.050.02'
$ws.Range("E48").Value = '  +0.05%  '
$ws.Range("E49").Value = '  -4.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05947'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.450'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.08%  '
